# The workbook originally has two sheets: "Sheet1" (the analysed data) and
# "Insights" (a sheet that only hosts a textbox/drawing summarising the
# dropped columns). This edit removes the "Insights" sheet (and its
# drawing/textbox) entirely, and renames "Sheet1" to "Analysed data(Task2)"
# so it becomes the sole, active/selected sheet in the workbook.

$wb = $excel.ActiveWorkbook

$excel.DisplayAlerts = $false

# Remove any shapes (the "Summary of Dropped Columns" textbox) living on the
# Insights sheet before removing the sheet itself.
$insights = $wb.Worksheets.Item("Insights")
foreach ($shp in $insights.Shapes) {
    $shp.Delete() | Out-Null
}
$insights.Delete() | Out-Null

# Rename the remaining data sheet.
$data = $wb.Worksheets.Item("Sheet1")
$data.Name = "Analysed data(Task2)"

# Make sure it is the active/selected sheet now that it is the only one.
$data.Select() | Out-Null

$excel.DisplayAlerts = $true
